$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "(WIT" + _GoBack bookmark + ") " -> merge into a single run "(WIT) "
#    and drop the old _GoBack bookmark (it gets relocated later).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Content.Find.Execute("(WIT) research group", $true, $false, $false, $false, $false, $true, 1, $false, "(WIT) research group", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Rewrite the "Although I am open to a variety..." paragraph.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Brian Caffo, Ciprian M. Crainiceanu, and Jeff Leek. After working with them for a few months, I see", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Brian Caffo (time series analysis), Ciprian M. Crainiceanu (WIT), and Scott L. Zeger (Bayesian statistics). After reading several papers in each of these groups, working with them for a few months, I see", `
    2) | Out-Null

# Place the _GoBack bookmark right after "...a vari" / before "ety of research"
$rng = $d.Content
$rng.Find.Execute("Although I am open to a vari", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pt = $rng.Duplicate
$pt.Collapse(0)
$d.Bookmarks.Add("_GoBack", $pt) | Out-Null

# ------------------------------------------------------------------
# 3) Track the deletion of "working with them for a few months, "
#    as a tracked change authored by Luchao Qi.
# ------------------------------------------------------------------
$word.Application.UserName = "Luchao Qi"
$d.TrackRevisions = $true

$delRng = $d.Content
$delRng.Find.Execute("working with them for a few months, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$delRng.Delete()

$d.TrackRevisions = $false
